$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "Prajwal_Niraj_Lawankar_BACV.docx"
$ws.Cells.Item(2, 2).Value = 33.33
$ws.Cells.Item(2, 3).Value = 43.45000076293945
$ws.Cells.Item(2, 4).Value = 1

$ws.Cells.Item(3, 1).Value = "Prajwal_Niraj_Lawankar_2025CV_1.pdf"
$ws.Cells.Item(3, 2).Value = 33.33
$ws.Cells.Item(3, 3).Value = 43.45000076293945
$ws.Cells.Item(3, 4).Value = 1

$ws.Cells.Item(4, 1).Value = "Prajwal_Lawankar_STC_CV_2025.pdf"
$ws.Cells.Item(4, 2).Value = 33.33
$ws.Cells.Item(4, 3).Value = 37.34000015258789
$ws.Cells.Item(4, 4).Value = 3

$ws.Cells.Item(5, 1).Value = "Mubashira_Khan_Operations_Admin_CV.docx21.pdf"
$ws.Cells.Item(5, 2).Value = 33.33
$ws.Cells.Item(5, 3).Value = 36.34000015258789
$ws.Cells.Item(5, 4).Value = 4

$ws.Cells.Item(6, 1).Value = "Prajwal_Lawankar_SC_CV_2025.pdf"
$ws.Cells.Item(6, 2).Value = 33.33
$ws.Cells.Item(6, 3).Value = 36.2400016784668
$ws.Cells.Item(6, 4).Value = 5

$ws.Cells.Item(7, 1).Value = "Prajwal_Lawankar___cV.pdf"
$ws.Cells.Item(7, 2).Value = 33.33
$ws.Cells.Item(7, 3).Value = 35.86999893188477
$ws.Cells.Item(7, 4).Value = 6

$ws.Cells.Item(8, 1).Value = "Prajwal_Niraj_Lawankar_BAE_CV.docx"
$ws.Cells.Item(8, 2).Value = 33.33
$ws.Cells.Item(8, 3).Value = 35.0099983215332
$ws.Cells.Item(8, 4).Value = 7

$ws.Cells.Item(9, 1).Value = "resume_prajwal_lawankar_11.pdf"
$ws.Cells.Item(9, 2).Value = 33.33
$ws.Cells.Item(9, 3).Value = 31.54999923706055
$ws.Cells.Item(9, 4).Value = 8

$ws.Cells.Item(10, 1).Value = "resume_prajwal_lawankar_1.pdf"
$ws.Cells.Item(10, 2).Value = 33.33
$ws.Cells.Item(10, 3).Value = 31.54999923706055
$ws.Cells.Item(10, 4).Value = 8

$ws.Cells.Item(11, 1).Value = "Prajwal_Niraj_Lawankar_CV.docx"
$ws.Cells.Item(11, 2).Value = 33.33
$ws.Cells.Item(11, 3).Value = 30.59000015258789
$ws.Cells.Item(11, 4).Value = 10

$ws.Cells.Item(12, 1).Value = "Shubham_UK_resume_for_reference.docx"
$ws.Cells.Item(12, 2).Value = 33.33
$ws.Cells.Item(12, 3).Value = 29.43000030517578
$ws.Cells.Item(12, 4).Value = 11

$ws.Cells.Item(13, 1).Value = "Omkar-resume-2024-DE2.pdf"
$ws.Cells.Item(13, 2).Value = 33.33
$ws.Cells.Item(13, 3).Value = 28.01000022888184
$ws.Cells.Item(13, 4).Value = 12

$ws.Cells.Item(14, 1).Value = "Phd_Omkar_CV.pdf"
$ws.Cells.Item(14, 2).Value = 33.33
$ws.Cells.Item(14, 3).Value = 26.71999931335449
$ws.Cells.Item(14, 4).Value = 13

$ws.Cells.Item(15, 1).Value = "Aishwarya_cv_for_reference.pdf"
$ws.Cells.Item(15, 2).Value = 33.33
$ws.Cells.Item(15, 3).Value = 21.21999931335449
$ws.Cells.Item(15, 4).Value = 14

$ws.Cells.Item(16, 1).Value = "PrajwalLawankar_GPM_CL.pdf"
$ws.Cells.Item(16, 2).Value = 0
$ws.Cells.Item(16, 3).Value = 39
$ws.Cells.Item(16, 4).Value = 15

$ws.Cells.Item(17, 1).Value = "Prajwal_Niraj_Lawankar_FSP_CV.docx"
$ws.Cells.Item(17, 2).Value = 0
$ws.Cells.Item(17, 3).Value = 33.81000137329102
$ws.Cells.Item(17, 4).Value = 16

$ws.Cells.Item(18, 1).Value = "Prajwal_Niraj_Lawankar_UK_CV.pdf"
$ws.Cells.Item(18, 2).Value = 0
$ws.Cells.Item(18, 3).Value = 29.3700008392334
$ws.Cells.Item(18, 4).Value = 17

$ws.Cells.Item(19, 1).Value = "Prajwal_Niraj_Lawankar_GS_CV.pdf"
$ws.Cells.Item(19, 2).Value = 0
$ws.Cells.Item(19, 3).Value = 29.28000068664551
$ws.Cells.Item(19, 4).Value = 18

$ws.Cells.Item(20, 1).Value = "Prajwal_Lawankar_UK_CV_Tailored.docx"
$ws.Cells.Item(20, 2).Value = 0
$ws.Cells.Item(20, 3).Value = 28.52000045776367
$ws.Cells.Item(20, 4).Value = 19

$ws.Cells.Item(21, 1).Value = "Winning_CV_Template_-_Extra_Curriculars.docx"
$ws.Cells.Item(21, 2).Value = 0
$ws.Cells.Item(21, 3).Value = 27.32999992370605
$ws.Cells.Item(21, 4).Value = 20

$ws.Cells.Item(22, 1).Value = "Prajwal_Niraj_Lawankar_Avaiva_CV.docx"
$ws.Cells.Item(22, 2).Value = 0
$ws.Cells.Item(22, 3).Value = 25.09000015258789
$ws.Cells.Item(22, 4).Value = 21

$ws.Cells.Item(23, 1).Value = "Prajwal_Niraj_Lawankar_Tatasteel_CV.docx"
$ws.Cells.Item(23, 2).Value = 0
$ws.Cells.Item(23, 3).Value = 23.79000091552734
$ws.Cells.Item(23, 4).Value = 22

$ws.Cells.Item(24, 1).Value = "Imerys_CV.pdf"
$ws.Cells.Item(24, 2).Value = 0
$ws.Cells.Item(24, 3).Value = 11.4399995803833
$ws.Cells.Item(24, 4).Value = 23

$ws.Cells.Item(25, 1).Value = "Tauheed_Kidwai_2.pdf"
$ws.Cells.Item(25, 2).Value = 0
$ws.Cells.Item(25, 3).Value = 10.69999980926514
$ws.Cells.Item(25, 4).Value = 24

